$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains one daily price record per row for "Espinaca" at
# "Vega Central Mapocho de Santiago". A new daily record needs to be
# inserted at row 270 (pushing the existing row 270..359 down by one to
# 271..360, matching the canonical diff's row-shift pattern and the new
# dimension A1:R360).
$ws.Rows("270:270").Insert()

# Populate the newly inserted row 270 with the new record's data.
$ws.Range("A270").Value = 9
$ws.Range("B270").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44588
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112012
$ws.Range("G270").Value = "Espinaca"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 97
$ws.Range("K270").Value = 10000
$ws.Range("L270").Value = 12000
$ws.Range("M270").Value = 10990
$ws.Range("N270").Value = "$/cuna 10 kilos"
$ws.Range("O270").Value = "Provincia de Chacabuco"
$ws.Range("P270").Value = 1099
$ws.Range("Q270").Value = 10
$ws.Range("R270").Value = "Hortaliza"
